$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": rows 3-7 move to "Ready for handoff" with a new
# handoff-generation timestamp (row 2 / 32d24a4a... stays "In Translation").
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$newDateOverview = "2017-02-09 17:15:39"
for ($r = 3; $r -le 7; $r++) {
    $wsOverview.Range("E$r").Value = "Ready for handoff"
    $wsOverview.Range("F$r").Value = "Ready for handoff"
    $wsOverview.Range("G$r").Value = $newDateOverview
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn": every row's Status becomes "Ready for handoff"; rows 3-7
# (the files that just got handed off) flip Priority "ht" -> "mt" and bump
# the Latest Handoff Datetime.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$newDateZhCn = "2017-02-09 17:15:19"
for ($r = 2; $r -le 7; $r++) {
    $wsZhCn.Range("C$r").Value = "Ready for handoff"
}
for ($r = 3; $r -le 7; $r++) {
    $wsZhCn.Range("E$r").Value = "mt"
    $wsZhCn.Range("H$r").Value = $newDateZhCn
}

# ---------------------------------------------------------------------------
# Sheet "de-de": same pattern as zh-cn, but with the de-de handoff timestamp.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$newDateDeDe = "2017-02-09 17:15:39"
for ($r = 2; $r -le 7; $r++) {
    $wsDeDe.Range("C$r").Value = "Ready for handoff"
}
for ($r = 3; $r -le 7; $r++) {
    $wsDeDe.Range("E$r").Value = "mt"
    $wsDeDe.Range("H$r").Value = $newDateDeDe
}
